$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-wide font change: Aptos Narrow -> Arial
$wb.Styles("Normal").Font.Name = "Arial"
$wb.Styles("Hyperlink").Font.Name = "Arial"

# New row of data (row 3)
$ws.Range("A3").Value = "يوسف يحيى محمد فاروق"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:faroukyoussif58@gmail.com")
$ws.Range("B3").Value = "faroukyoussif58@gmail.com"
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("C3").Value = "https://github.com/YoussifYahia/Security-Task.git"

$ws.Range("C3").Select() | Out-Null
